$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# 1) Intro paragraph: "...how easy it is read your book..." ->
#    "...how easy it is to read your book..." with the inserted "to " word
#    split into its own run, and the "_GoBack" bookmark (which used to sit
#    at the very end of the document) moved here, collapsed between
#    "to " and "read your book...".
# -------------------------------------------------------------------------
$introPara = $d.Paragraphs.Item(2)
$introText = $introPara.Range.Text
$readIdx = $introText.IndexOf("read your book")
$splitPos = $introPara.Range.Start + $readIdx

# Insert the missing word "to " right before "read your book...".
$insertRng = $d.Range($splitPos, $splitPos)
$insertRng.InsertBefore("to ")

# Re-locate the boundary between "...is " and "to read your book..." (this
# is where the first run needs to be cut so "to " becomes its own run) and
# the boundary between "to " and "read your book..." (this is where the
# "_GoBack" bookmark belongs).
$introText2 = $introPara.Range.Text
$toIdx = $introText2.IndexOf("to read your book")
$beforeToPos = $introPara.Range.Start + $toIdx
$afterToPos = $beforeToPos + 3

# Use a throwaway bookmark to force a clean run split before "to " (no
# leftover rPr/formatting marks get left behind); it is removed again once
# the real bookmark has been added after it.
$tempRng = $d.Range($beforeToPos, $beforeToPos)
$d.Bookmarks.Add("zzzTempSplit", $tempRng)

# Bookmark names are unique, so adding "_GoBack" here automatically moves
# it away from wherever it used to be (the end of the document).
$goBackRng = $d.Range($afterToPos, $afterToPos)
$d.Bookmarks.Add("_GoBack", $goBackRng)

$d.Bookmarks("zzzTempSplit").Delete()

# -------------------------------------------------------------------------
# 2) "Making Text Larger" heading: add a <w:lastRenderedPageBreak/> marker
#    before the run's text.
# -------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Trim() -eq "Making Text Larger") {
        $xmlFrag = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' w:rsidR='004E21B4' w:rsidRDefault='004E21B4' w:rsidP='004E21B4'><w:pPr><w:pStyle w:val='Heading1'/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Making Text Larger</w:t></w:r></w:p>"
        $para.Range.InsertXML($xmlFrag)
        break
    }
}
